$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the CreatedAt timestamp in A1
$ws.Range("A1").Value = "CreatedAt: 2026-02-13T18:07:24"

# Update intertie pricing values (columns V-Z) across data rows
$ws.Range("V4").Value = 128.58
$ws.Range("W4").Value = 124.62
$ws.Range("X4").Value = 57.87
$ws.Range("Y4").Value = 60.06
$ws.Range("Z4").Value = 55.88
$ws.Range("V6").Value = -4.76
$ws.Range("W6").Value = -4.86
$ws.Range("X6").Value = -2.03
$ws.Range("Y6").Value = -1.38
$ws.Range("Z6").Value = -0.78
$ws.Range("V9").Value = 126.75
$ws.Range("W9").Value = 123.55
$ws.Range("X9").Value = 57.21
$ws.Range("Y9").Value = 60.12
$ws.Range("Z9").Value = 56.04
$ws.Range("V11").Value = -6.59
$ws.Range("W11").Value = -5.93
$ws.Range("X11").Value = -2.69
$ws.Range("Y11").Value = -1.32
$ws.Range("V14").Value = 126.75
$ws.Range("W14").Value = 123.55
$ws.Range("X14").Value = 57.21
$ws.Range("Y14").Value = 60.12
$ws.Range("Z14").Value = 56.04
$ws.Range("V16").Value = -6.59
$ws.Range("W16").Value = -5.93
$ws.Range("X16").Value = -2.69
$ws.Range("Y16").Value = -1.32
$ws.Range("V19").Value = 54.48
$ws.Range("W19").Value = 61.9
$ws.Range("X19").Value = 46.5
$ws.Range("Y19").Value = 59.2
$ws.Range("Z19").Value = 55.17
$ws.Range("V20").Value = -72.27
$ws.Range("W20").Value = -60.95
$ws.Range("X20").Value = -10.54
$ws.Range("Y20").Value = 0
$ws.Range("V21").Value = -6.59
$ws.Range("W21").Value = -6.63
$ws.Range("X21").Value = -2.85
$ws.Range("Y21").Value = -2.25
$ws.Range("Z21").Value = -1.49
$ws.Range("V24").Value = 126.75
$ws.Range("W24").Value = 122.85
$ws.Range("X24").Value = 57.04
$ws.Range("Y24").Value = 59.2
$ws.Range("Z24").Value = 55.17
$ws.Range("V26").Value = -6.59
$ws.Range("W26").Value = -6.63
$ws.Range("X26").Value = -2.85
$ws.Range("Y26").Value = -2.25
$ws.Range("Z26").Value = -1.49
$ws.Range("V29").Value = 52.11
$ws.Range("W29").Value = 59.72
$ws.Range("X29").Value = 45.44
$ws.Range("Y29").Value = 58.02
$ws.Range("Z29").Value = 54.22
$ws.Range("V30").Value = -72.27
$ws.Range("W30").Value = -60.95
$ws.Range("X30").Value = -10.54
$ws.Range("Y30").Value = 0
$ws.Range("V31").Value = -8.960000000000001
$ws.Range("W31").Value = -8.81
$ws.Range("X31").Value = -3.92
$ws.Range("Y31").Value = -3.42
$ws.Range("Z31").Value = -2.44
$ws.Range("V34").Value = 128.21
$ws.Range("W34").Value = 125.34
$ws.Range("X34").Value = 57.87
$ws.Range("Y34").Value = 61.2
$ws.Range("Z34").Value = 57.12
$ws.Range("V36").Value = -5.13
$ws.Range("W36").Value = -4.14
$ws.Range("X36").Value = -2.03
$ws.Range("Y36").Value = -0.24
$ws.Range("Z36").Value = 0.46
$ws.Range("V39").Value = 128.58
$ws.Range("W39").Value = 124.62
$ws.Range("X39").Value = 57.87
$ws.Range("Y39").Value = 60.06
$ws.Range("Z39").Value = 55.88
$ws.Range("V41").Value = -4.76
$ws.Range("W41").Value = -4.86
$ws.Range("X41").Value = -2.03
$ws.Range("Y41").Value = -1.38
$ws.Range("Z41").Value = -0.78
$ws.Range("V44").Value = 134.68
$ws.Range("W44").Value = 130.79
$ws.Range("X44").Value = 60.5
$ws.Range("Y44").Value = 62.25
$ws.Range("Z44").Value = 57.64
$ws.Range("V46").Value = 1.35
$ws.Range("W46").Value = 1.31
$ws.Range("X46").Value = 0.61
$ws.Range("Y46").Value = 0.8100000000000001
$ws.Range("Z46").Value = 0.98
$ws.Range("V49").Value = 130.34
$ws.Range("W49").Value = 134.6
$ws.Range("X49").Value = 62.13
$ws.Range("Y49").Value = 63.28
$ws.Range("Z49").Value = 58.23
$ws.Range("V51").Value = -3
$ws.Range("W51").Value = 5.11
$ws.Range("X51").Value = 2.24
$ws.Range("Y51").Value = 1.84
$ws.Range("Z51").Value = 1.57
$ws.Range("V54").Value = 130.09
$ws.Range("W54").Value = 126.2
$ws.Range("X54").Value = 58.61
$ws.Range("Y54").Value = 61.2
$ws.Range("Z54").Value = 56.55
$ws.Range("V56").Value = -3.25
$ws.Range("W56").Value = -3.28
$ws.Range("X56").Value = -1.29
$ws.Range("Y56").Value = -0.24
$ws.Range("Z56").Value = -0.11
$ws.Range("V59").Value = 138.6
$ws.Range("W59").Value = 134.74
$ws.Range("X59").Value = 62.33
$ws.Range("Y59").Value = 64.01000000000001
$ws.Range("Z59").Value = 59.08
$ws.Range("V61").Value = 5.27
$ws.Range("W61").Value = 5.25
$ws.Range("X61").Value = 2.43
$ws.Range("Y61").Value = 2.56
$ws.Range("Z61").Value = 2.42
$ws.Range("V64").Value = 140.95
$ws.Range("W64").Value = 136.73
$ws.Range("X64").Value = 63.45
$ws.Range("Y64").Value = 65.09
$ws.Range("Z64").Value = 59.96
$ws.Range("V66").Value = 7.61
$ws.Range("W66").Value = 7.25
$ws.Range("X66").Value = 3.55
$ws.Range("Y66").Value = 3.65
$ws.Range("Z66").Value = 3.3
$ws.Range("W69").Value = 138.33
$ws.Range("X69").Value = 63.99
$ws.Range("Y69").Value = 65.79000000000001
$ws.Range("Z69").Value = 60.6
$ws.Range("V71").Value = 8.66
$ws.Range("W71").Value = 8.85
$ws.Range("X71").Value = 4.1
$ws.Range("Y71").Value = 4.34
$ws.Range("Z71").Value = 3.94
$ws.Range("V74").Value = 138.17
$ws.Range("W74").Value = 134.6
$ws.Range("X74").Value = 62.26
$ws.Range("Y74").Value = 64.01000000000001
$ws.Range("Z74").Value = 58.96
$ws.Range("V76").Value = 4.84
$ws.Range("W76").Value = 5.11
$ws.Range("X76").Value = 2.37
$ws.Range("Y76").Value = 2.56
$ws.Range("Z76").Value = 2.3
$ws.Range("V79").Value = 133.34
$ws.Range("W79").Value = 129.48
$ws.Range("X79").Value = 59.9
$ws.Range("Y79").Value = 61.45
$ws.Range("Z79").Value = 56.66
$ws.Range("V84").Value = 124.97
$ws.Range("W84").Value = 122.96
$ws.Range("X84").Value = 56.99
$ws.Range("Y84").Value = 61.45
$ws.Range("Z84").Value = 56.94
$ws.Range("V86").Value = -8.369999999999999
$ws.Range("W86").Value = -6.52
$ws.Range("X86").Value = -2.91
$ws.Range("Y86").Value = 0
$ws.Range("Z86").Value = 0.28
$ws.Range("V89").Value = 124.38
$ws.Range("W89").Value = 120.67
$ws.Range("X89").Value = 55.98
$ws.Range("Y89").Value = 58.02
$ws.Range("Z89").Value = 54.17
$ws.Range("V91").Value = -8.960000000000001
$ws.Range("W91").Value = -8.81
$ws.Range("X91").Value = -3.92
$ws.Range("Y91").Value = -3.42
$ws.Range("Z91").Value = -2.49
